# Scheduled-runner price/profit refresh for the Durandal_Profits workbook.
# Updates the market-price-derived columns (H:N) on a handful of leve rows
# across several crafting-job sheets. Columns:
#   H  currentAveragePrice
#   I  currentAveragePriceNQ
#   J  currentAveragePriceHQ
#   K  LevePriceNQ
#   L  LevePriceHQ
#   M  LeveProfitNQ
#   N  LeveProfitHQ

$wb = $excel.ActiveWorkbook

# --- ALC ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H18").Value = 1100
$ws.Range("I18").Value = 1000
$ws.Range("J18").Value = 1150
$ws.Range("K18").Value = 1000
$ws.Range("L18").Value = 1150
$ws.Range("M18").Value = -716
$ws.Range("N18").Value = -1718

$ws.Range("H62").Value = 1641.9048
$ws.Range("I62").Value = 1583.3334
$ws.Range("J62").Value = 1993.3334
$ws.Range("K62").Value = 1583.3334
$ws.Range("L62").Value = 1993.3334
$ws.Range("M62").Value = -959.3334
$ws.Range("N62").Value = -3241.3334

$ws.Range("H65").Value = 1641.9048
$ws.Range("I65").Value = 1583.3334
$ws.Range("J65").Value = 1993.3334
$ws.Range("K65").Value = 7916.666999999999
$ws.Range("L65").Value = 9966.666999999999
$ws.Range("M65").Value = -4796.666999999999
$ws.Range("N65").Value = -16206.667

$ws.Range("H98").Value = 5343831.5
$ws.Range("I98").Value = 58333.95
$ws.Range("J98").Value = 55556056
$ws.Range("K98").Value = 58333.95
$ws.Range("L98").Value = 55556056
$ws.Range("M98").Value = -56835.95
$ws.Range("N98").Value = -55559052

$ws.Range("H122").Value = 5343831.5
$ws.Range("I122").Value = 58333.95
$ws.Range("J122").Value = 55556056
$ws.Range("K122").Value = 175001.85
$ws.Range("L122").Value = 166668168
$ws.Range("M122").Value = -172551.85
$ws.Range("N122").Value = -166673068

$ws.Range("H141").Value = 1680.75
$ws.Range("I141").Value = 1234
$ws.Range("J141").Value = 3021
$ws.Range("K141").Value = 3702
$ws.Range("L141").Value = 9063
$ws.Range("M141").Value = 1478
$ws.Range("N141").Value = -19423

# --- ARM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 351052.62
$ws.Range("I32").Value = 4990.305
$ws.Range("J32").Value = 1201789.2
$ws.Range("K32").Value = 4990.305
$ws.Range("L32").Value = 1201789.2
$ws.Range("M32").Value = -4703.305
$ws.Range("N32").Value = -1202363.2

$ws.Range("H37").Value = 166672600
$ws.Range("I37").Value = 500001020
$ws.Range("J37").Value = 8399
$ws.Range("K37").Value = 500001020
$ws.Range("L37").Value = 8399
$ws.Range("M37").Value = -500000747
$ws.Range("N37").Value = -8945

$ws.Range("H74").Value = 1244
$ws.Range("I74").Value = 687.6
$ws.Range("J74").Value = 2171.3333
$ws.Range("K74").Value = 687.6
$ws.Range("L74").Value = 2171.3333
$ws.Range("M74").Value = 186.4
$ws.Range("N74").Value = -3919.3333

$ws.Range("H77").Value = 1244
$ws.Range("I77").Value = 687.6
$ws.Range("J77").Value = 2171.3333
$ws.Range("K77").Value = 3438
$ws.Range("L77").Value = 10856.6665
$ws.Range("M77").Value = 930
$ws.Range("N77").Value = -19592.6665

$ws.Range("H132").Value = 1672.9108
$ws.Range("I132").Value = 1574.2821
$ws.Range("K132").Value = 4722.846299999999
$ws.Range("M132").Value = -2192.846299999999

# --- BSM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H94").Value = 1120.0869
$ws.Range("I94").Value = 899.7778
$ws.Range("J94").Value = 1913.2
$ws.Range("K94").Value = 899.7778
$ws.Range("L94").Value = 1913.2
$ws.Range("M94").Value = -448.7778
$ws.Range("N94").Value = -2815.2

$ws.Range("H140").Value = 69296.25
$ws.Range("J140").Value = 69296.25
$ws.Range("L140").Value = 69296.25
$ws.Range("N140").Value = -79656.25

# --- CRP -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H50").Value = 8264.111000000001
$ws.Range("J50").Value = 8264.111000000001
$ws.Range("L50").Value = 8264.111000000001
$ws.Range("N50").Value = -9514.111000000001

$ws.Range("H60").Value = 8474.5
$ws.Range("J60").Value = 8474.5
$ws.Range("L60").Value = 8474.5
$ws.Range("N60").Value = -9496.5

$ws.Range("H74").Value = 14663.7
$ws.Range("J74").Value = 16039.111
$ws.Range("L74").Value = 16039.111
$ws.Range("N74").Value = -17787.111

$ws.Range("H77").Value = 14663.7
$ws.Range("J77").Value = 16039.111
$ws.Range("L77").Value = 48117.333
$ws.Range("N77").Value = -56853.333

$ws.Range("H138").Value = 44250
$ws.Range("J138").Value = 44250
$ws.Range("L138").Value = 44250
$ws.Range("N138").Value = -54530

$ws.Range("H140").Value = 89340
$ws.Range("J140").Value = 89340
$ws.Range("L140").Value = 89340
$ws.Range("N140").Value = -99700

# --- CUL -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H131").Value = 794.49
$ws.Range("I131").Value = 403.07693
$ws.Range("J131").Value = 852.977
$ws.Range("K131").Value = 1209.23079
$ws.Range("L131").Value = 2558.931
$ws.Range("M131").Value = 3830.76921
$ws.Range("N131").Value = -12638.931

$ws.Range("H140").Value = 3718.7693
$ws.Range("I140").Value = 2333.2727
$ws.Range("J140").Value = 5511.7646
$ws.Range("K140").Value = 6999.8181
$ws.Range("L140").Value = 16535.2938
$ws.Range("M140").Value = -1819.8181
$ws.Range("N140").Value = -26895.2938

# --- GSM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H2").Value = 38.133335
$ws.Range("I2").Value = 17.2
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 17.2
$ws.Range("L2").Value = 80
$ws.Range("M2").Value = 95.8
$ws.Range("N2").Value = -306

$ws.Range("H122").Value = 2995.6667
$ws.Range("I122").Value = 2995.6667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8987.000100000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6537.000100000001
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 2284.375
$ws.Range("I132").Value = 1819.4193
$ws.Range("J132").Value = 3885.889
$ws.Range("K132").Value = 5458.257900000001
$ws.Range("L132").Value = 11657.667
$ws.Range("M132").Value = -2928.257900000001
$ws.Range("N132").Value = -16717.667

$ws.Range("H140").Value = 99853
$ws.Range("J140").Value = 99853
$ws.Range("L140").Value = 99853
$ws.Range("N140").Value = -110213

# --- WVR -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H132").Value = 23810772
$ws.Range("I132").Value = 33334368
$ws.Range("J132").Value = 1782.2778
$ws.Range("K132").Value = 100003104
$ws.Range("L132").Value = 5346.8334
$ws.Range("M132").Value = -100000574
$ws.Range("N132").Value = -10406.8334

$ws.Range("H136").Value = 551.3871
$ws.Range("I136").Value = 471.3654
$ws.Range("J136").Value = 967.5
$ws.Range("K136").Value = 1414.0962
$ws.Range("L136").Value = 2902.5
$ws.Range("M136").Value = 1135.9038
$ws.Range("N136").Value = -8002.5
